$wb = $excel.ActiveWorkbook

# Locate the "QuickLink" worksheet (falls back to the active sheet, which is
# this sheet in the source workbook, in case the name already changed).
$ws = $null
try {
    $ws = $wb.Worksheets.Item("QuickLink")
} catch {
    $ws = $null
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

# Rename the "QuickLink" sheet to "QuickLinks"
$ws.Name = "QuickLinks"

# Update the sheet's title cell (A1) to match the new sheet name
$ws.Range("A1").Value = "Quick Links"

# Reset the lingering stale selection (previously parked at A17) back to A1
$ws.Range("A1").Select()

Write-Host "Renamed QuickLink sheet to QuickLinks and updated title cell"
